$wb = $excel.ActiveWorkbook

$rows = @(
    @{ Sheet = "ROW35-FE-LIFTER";  A = "2025-03-07 05:42:06"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x14,0x41,0x0c,"; E = "0x d"; G = "568631262647113770877196"; I = 13 },
    @{ Sheet = "ROW35-MID-LIFTER"; A = "2025-03-07 05:29:35"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x15,0x41,0x0c,"; E = "0x e"; G = "568631262647113770942732"; I = 14 },
    @{ Sheet = "ROW02-FE-LIFTER";  A = "2025-03-07 05:51:45"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x06,0x41,0x0c,"; E = "0xff"; G = "568631262647113769959692"; I = 255 },
    @{ Sheet = "ROW02-MID-LIFTER"; A = "2025-03-07 05:41:15"; C = "0x78,0x69,0x90,0x01,0x00,0x00,0x18,0x0b,0x40,0x0c,"; E = "0x 3"; G = "568631262647113769959692"; I = 3 }
)

foreach ($r in $rows) {
    $ws = $wb.Worksheets.Item($r.Sheet)
    $row = 70

    $ws.Cells.Item($row, 1).Value = $r.A
    $ws.Cells.Item($row, 2).Value = "0x01,0x90 "
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = "0x01,0x90,"
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Value = 400

    # Column G holds a long digit-only string (too large to round-trip as a
    # double without losing precision), so force text formatting before
    # assigning it - otherwise it gets auto-coerced into a number.
    $ws.Cells.Item($row, 7).NumberFormat = "@"
    $ws.Cells.Item($row, 7).Value = $r.G

    $ws.Cells.Item($row, 8).Value = 400
    $ws.Cells.Item($row, 9).Value = $r.I
}

$wb.Save()
